# Weekly cryptos-list refresh (commit message: "Updated cryptos list on Wed Mar 27
# 19:16:33 UTC 2024 with GitHub Actions"). Updates the Price (D) and Volume(1h) (E)
# columns scraped for each coin; rows 48/49 additionally swap places (Stellar now
# ranks above THORChain), so Coin (B) and Link (C) are rewritten there too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column cells are stored as TEXT in the source sheet (e.g. "68.734.79" / "1.00"),
# not numbers. A leading apostrophe forces Excel to keep numeric-looking entries as text
# instead of silently converting them to a Number (which would drop the trailing zero in
# "1.00" or garble "68.734.79" as it is not a valid numeric literal anyway).

$ws.Range('D2').Value = '68.734.79'
$ws.Range('E2').Value = '  -1.64%  '

$ws.Range('D3').Value = '3.480.81'
$ws.Range('E3').Value = '  -2.47%  '

$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').Value = '''567.69'
$ws.Range('E5').Value = '  -1.88%  '

$ws.Range('D6').Value = '''182.23'
$ws.Range('E6').Value = '  -3.32%  '

$ws.Range('E7').Value = '  -3.12%  '

$ws.Range('D8').Value = '3.468.07'
$ws.Range('E8').Value = '  -2.67%  '

$ws.Range('E9').Value = '  +0.17%  '

$ws.Range('E10').Value = '  +3.52%  '

$ws.Range('D11').Value = '''0.640'
$ws.Range('E11').Value = '  -2.86%  '

$ws.Range('D12').Value = '''53.62'
$ws.Range('E12').Value = '  -3.93%  '

$ws.Range('E13').Value = '  -0.98%  '

$ws.Range('D14').Value = '''9.38'
$ws.Range('E14').Value = '  -1.93%  '

$ws.Range('D15').Value = '4.033.84'
$ws.Range('E15').Value = '  -2.69%  '

$ws.Range('D16').Value = '''19.12'
$ws.Range('E16').Value = '  -3.38%  '

$ws.Range('D17').Value = '68.595.93'
$ws.Range('E17').Value = '  -1.57%  '

$ws.Range('D18').Value = '3.467.44'
$ws.Range('E18').Value = '  -3.00%  '

$ws.Range('D19').Value = '''12.21'
$ws.Range('E19').Value = '  -3.11%  '

$ws.Range('E20').Value = '  -1.35%  '

$ws.Range('D21').Value = '''537.07'
$ws.Range('E21').Value = '  +13.14%  '

$ws.Range('E22').Value = '  -2.94%  '

$ws.Range('D23').Value = '''19.04'
$ws.Range('E23').Value = '  -1.00%  '

$ws.Range('D24').Value = '''4.94'
$ws.Range('E24').Value = '  -2.54%  '

$ws.Range('D25').Value = '''4.37'
$ws.Range('E25').Value = '  +0.46%  '

$ws.Range('D26').Value = '''93.90'
$ws.Range('E26').Value = '  +0.56%  '

$ws.Range('D27').Value = '''2.89'
$ws.Range('E27').Value = '  -4.21%  '

$ws.Range('D28').Value = '''10.70'
$ws.Range('E28').Value = '  -2.54%  '

$ws.Range('D29').Value = '''9.00'
$ws.Range('E29').Value = '  -3.27%  '

$ws.Range('D30').Value = '''31.33'
$ws.Range('E30').Value = '  -2.75%  '

$ws.Range('D31').Value = '''7.12'
$ws.Range('E31').Value = '  -7.99%  '

$ws.Range('D32').Value = '''12.45'
$ws.Range('E32').Value = '  +2.25%  '

$ws.Range('D33').Value = '''64.34'
$ws.Range('E33').Value = '  -2.80%  '

$ws.Range('E34').Value = '  -5.28%  '

$ws.Range('D35').Value = '''568.12'
$ws.Range('E35').Value = '  -1.98%  '

$ws.Range('D36').Value = '''1.00'
$ws.Range('E36').Value = '  +0.11%  '

$ws.Range('D37').Value = '''37.56'

$ws.Range('E38').Value = '  -0.44%  '

$ws.Range('D39').Value = '''2.98'
$ws.Range('E39').Value = '  +4.55%  '

$ws.Range('E40').Value = '  -4.66%  '

$ws.Range('E41').Value = '  -4.95%  '

$ws.Range('D42').Value = '''0.131'
$ws.Range('E42').Value = '  -5.55%  '

$ws.Range('E43').Value = '  -4.76%  '

$ws.Range('D44').Value = '3.219.60'
$ws.Range('E44').Value = '  -0.19%  '

$ws.Range('D45').Value = '''2.94'
$ws.Range('E45').Value = '  -4.38%  '

$ws.Range('D46').Value = '''3.43'
$ws.Range('E46').Value = '  +2.05%  '

$ws.Range('E47').Value = '  -1.53%  '

$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').Value = '''0.134'
$ws.Range('E48').Value = '  -2.98%  '

$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').Value = '''8.93'
$ws.Range('E49').Value = '  -5.87%  '

$ws.Range('D50').Value = '''0.998'
$ws.Range('E50').Value = '  -0.26%  '

$ws.Range('D51').Value = '''137.65'
$ws.Range('E51').Value = '  +1.70%  '
